# Update "想去人数" (number of people wanting to attend) figures that changed
# between the previous data pull and the one generated at commit 456a3b4.
#
# The same underlying events appear on both the "展览" sheet (sheet1) and the
# "全部类型" aggregate sheet (sheet4); both need to be kept in sync.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Cell -> new value updates for the "展览" sheet
$updates1 = @{
    "F3"  = 521
    "F4"  = 1075
    "F6"  = 24
    "F9"  = 1104
    "F10" = 15444
    "F11" = 222
    "F17" = 64
    "F19" = 104
    "F20" = 1261
    "F28" = 15
    "F30" = 136
    "F31" = 10957
    "F34" = 108
    "F36" = 3786
}

foreach ($cell in $updates1.Keys) {
    $ws1.Range($cell).Value = $updates1[$cell]
}

# Cell -> new value updates for the "全部类型" sheet (same events, offset rows)
$updates4 = @{
    "F4"  = 521
    "F5"  = 1075
    "F7"  = 24
    "F10" = 1104
    "F11" = 15444
    "F12" = 222
    "F18" = 64
    "F20" = 104
    "F21" = 1261
    "F29" = 15
    "F31" = 136
    "F33" = 10957
    "F36" = 108
    "F38" = 3786
}

foreach ($cell in $updates4.Keys) {
    $ws4.Range($cell).Value = $updates4[$cell]
}
